# Updated rmi files 3.4.3
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# "About" sheet: remove the stray date in C1; add two new footnote rows
# (A13/A14) with an explanatory note about US exemptions. The narrative
# text in A3:A11 (control-lever boilerplate + boolean-lever explanation)
# is unchanged.
# ---------------------------------------------------------------------
$wsAbout = $wb.Worksheets.Item("About")

$wsAbout.Range("C1").Clear()

$wsAbout.Range("A13").Value = "In the U.S., we exempt agriculture and water and waste process emissions. Generally, "
$wsAbout.Range("A14").Value = "proposed taxes do not cover these sectors."

# ---------------------------------------------------------------------
# "BEPEfCT" sheet: was a single boolean lever ("Boolean" / 0). It becomes
# a per-industry-sector table of booleans, one row per sector, with the
# label cell restyled from a date format to an italic "Unit" caption.
# ---------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("BEPEfCT")

$wsData.Range("A1").Value = "Unit: boolean (0 or 1)"
$wsData.Range("A1").Font.Italic = $true
$wsData.Range("B1").Value = "Exempt Process Emissions from Carbon Tax"

$sectors = @(
    @{ Name = "agriculture and forestry 01T03"; Value = 1 },
    @{ Name = "coal mining 05"; Value = 0 },
    @{ Name = "oil and gas extraction 06"; Value = 0 },
    @{ Name = "other mining and quarrying 07T08"; Value = 0 },
    @{ Name = "food beverage and tobacco 10T12"; Value = 0 },
    @{ Name = "textiles apparel and leather 13T15"; Value = 0 },
    @{ Name = "wood products 16"; Value = 0 },
    @{ Name = "pulp paper and printing 17T18"; Value = 0 },
    @{ Name = "refined petroleum and coke 19"; Value = 0 },
    @{ Name = "chemicals 20"; Value = 0 },
    @{ Name = "rubber and plastic products 22"; Value = 0 },
    @{ Name = "glass and glass products 231"; Value = 0 },
    @{ Name = "cement and other nonmetallic minerals 239"; Value = 0 },
    @{ Name = "iron and steel 241"; Value = 0 },
    @{ Name = "other metals 242"; Value = 0 },
    @{ Name = "metal products except machinery and vehicles 25"; Value = 0 },
    @{ Name = "computers and electronics 26"; Value = 0 },
    @{ Name = "appliances and electrical equipment 27"; Value = 0 },
    @{ Name = "other machinery 28"; Value = 0 },
    @{ Name = "road vehicles 29"; Value = 0 },
    @{ Name = "nonroad vehicles 30"; Value = 0 },
    @{ Name = "other manufacturing 31T33"; Value = 0 },
    @{ Name = "energy pipelines and gas processing 352T353"; Value = 0 },
    @{ Name = "water and waste 36T39"; Value = 1 },
    @{ Name = "construction 41T43"; Value = 0 }
)

$r = 2
foreach ($sector in $sectors) {
    $wsData.Cells.Item($r, 1).Value = $sector.Name
    $wsData.Cells.Item($r, 2).Value = $sector.Value
    $r = $r + 1
}

$wsData.Columns.Item(1).ColumnWidth = 47
$wsData.PageSetup.Orientation = 1
